$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-five character (U+2085) used in ShibaInu price (row 16)
$sub5 = [char]0x2085

# Row 2
$ws.Range("D2").Value = "'25.841.47"
$ws.Range("E2").Value = "'  +0.39%  "

# Row 3
$ws.Range("D3").Value = "'1.635.98"
$ws.Range("E3").Value = "'  +0.52%  "

# Row 4
$ws.Range("E4").Value = "'  +0.11%  "

# Row 5
$ws.Range("D5").Value = "'215.20"
$ws.Range("E5").Value = "'  -0.10%  "

# Row 6
$ws.Range("D6").Value = "'0.5087"
$ws.Range("E6").Value = "'  -0.39%  "

# Row 7
$ws.Range("E7").Value = "'  +0.11%  "

# Row 8
$ws.Range("D8").Value = "'0.2584"
$ws.Range("E8").Value = "'  +0.87%  "

# Row 9
$ws.Range("D9").Value = "'0.06432"
$ws.Range("E9").Value = "'  +1.80%  "

# Row 10
$ws.Range("D10").Value = "'20.36"
$ws.Range("E10").Value = "'  +4.90%  "

# Row 11
$ws.Range("D11").Value = "'0.07798"
$ws.Range("E11").Value = "'  +0.24%  "

# Row 12
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.665.89"
$ws.Range("E12").Value = "'  +2.36%  "

# Row 13
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.251"
$ws.Range("E13").Value = "'  +0.41%  "

# Row 14
$ws.Range("D14").Value = "'1.862.11"
$ws.Range("E14").Value = "'  +0.51%  "

# Row 15
$ws.Range("D15").Value = "'0.5597"
$ws.Range("E15").Value = "'  +1.71%  "

# Row 16
$v_D16 = "'0.0{0}7668" -f $sub5
$ws.Range("D16").Value = $v_D16
$ws.Range("E16").Value = "'  +1.75%  "

# Row 17
$ws.Range("E17").Value = "'  -0.47%  "

# Row 18
$ws.Range("D18").Value = "'25.849.42"
$ws.Range("E18").Value = "'  +0.22%  "

# Row 19
$ws.Range("E19").Value = "'  +0.15%  "

# Row 20
$ws.Range("D20").Value = "'4.376"
$ws.Range("E20").Value = "'  -0.59%  "

# Row 21
$ws.Range("D21").Value = "'192.56"
$ws.Range("E21").Value = "'  -0.71%  "

# Row 22
$ws.Range("D22").Value = "'9.943"
$ws.Range("E22").Value = "'  +1.22%  "

# Row 23
$ws.Range("E23").Value = "'  +2.53%  "

# Row 24
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "'  +0.02%  "

# Row 25
$ws.Range("D25").Value = "'1.746"
$ws.Range("E25").Value = "'  -7.51%  "

# Row 26
$ws.Range("D26").Value = "'139.21"
$ws.Range("E26").Value = "'  -2.01%  "

# Row 27
$ws.Range("D27").Value = "'0.1234"
$ws.Range("E27").Value = "'  -1.62%  "

# Row 28
$ws.Range("D28").Value = "'6.832"
$ws.Range("E28").Value = "'  +1.37%  "

# Row 29
$ws.Range("E29").Value = "'  -0.05%  "

# Row 30
$ws.Range("D30").Value = "'1.240"
$ws.Range("E30").Value = "'  +0.15%  "

# Row 31
$ws.Range("D31").Value = "'0.04968"
$ws.Range("E31").Value = "'  +1.92%  "

# Row 32
$ws.Range("D32").Value = "'3.308"
$ws.Range("E32").Value = "'  +2.53%  "

# Row 33
$ws.Range("D33").Value = "'3.259"
$ws.Range("E33").Value = "'  +2.83%  "

# Row 34
$ws.Range("D34").Value = "'1.570"
$ws.Range("E34").Value = "'  +1.99%  "

# Row 35
$ws.Range("E35").Value = "'  +0.57%  "

# Row 36
$ws.Range("D36").Value = "'0.9018"
$ws.Range("E36").Value = "'  +1.01%  "

# Row 37
$ws.Range("D37").Value = "'0.5575"
$ws.Range("E37").Value = "'  +1.36%  "

# Row 38
$ws.Range("D38").Value = "'2.568"
$ws.Range("E38").Value = "'  +1.29%  "

# Row 39
$ws.Range("D39").Value = "'1.134.00"
$ws.Range("E39").Value = "'  +1.93%  "

# Row 40
$ws.Range("D40").Value = "'0.01571"
$ws.Range("E40").Value = "'  +1.40%  "

# Row 41
$ws.Range("D41").Value = "'0.9966"
$ws.Range("E41").Value = "'  -0.49%  "

# Row 42
$ws.Range("D42").Value = "'99.15"
$ws.Range("E42").Value = "'  +2.02%  "

# Row 43
$ws.Range("D43").Value = "'5.459"
$ws.Range("E43").Value = "'  -1.89%  "

# Row 44
$ws.Range("D44").Value = "'0.8003"
$ws.Range("E44").Value = "'  +0.76%  "

# Row 45
$ws.Range("E45").Value = "'  +0.37%  "

# Row 46
$ws.Range("D46").Value = "'55.67"
$ws.Range("E46").Value = "'  +2.00%  "

# Row 47
$ws.Range("E47").Value = "'  -3.74%  "

# Row 48
$ws.Range("D48").Value = "'7.825"
$ws.Range("E48").Value = "'  +3.58%  "

# Row 49
$ws.Range("D49").Value = "'0.05026"
$ws.Range("E49").Value = "'  -2.10%  "

# Row 50
$ws.Range("D50").Value = "'0.9968"
$ws.Range("E50").Value = "'  -0.28%  "

# Row 51
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "'  +0.25%  "
